$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 ("I0") and J1 ("IF"), styled like the other headers (copy H1's format) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values for I2:J55 ---
$iValues = @(6,8,7,4,7,7,6,7,4,7,9,9,5,5,8,1,7,8,6,6,5,7,7,8,6,6,7,6,6,9,5,9,8,7,8,5,7,3,10,9,5,10,7,5,7,8,5,5,7,8,6,6,8,7)
$jValues = @(6,8,7,4,8,7,7,7,4,7,9,9,6,5,8,2,7,8,7,6,6,8,7,9,6,6,7,6,7,9,6,9,8,7,8,5,7,4,10,9,5,10,7,5,8,9,6,5,7,8,6,6,8,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
